# Apply crypto price/volume updates from the Coinranking scrape refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Note: several 'Price' values in column D look like plain numbers (e.g. 1.014,
# 88.00, 0.5160) but must stay TEXT, exactly as scraped, including trailing
# zeros. We force text by prefixing the literal with an apostrophe, the same
# trick Excel's UI uses (sets the cell's quote-prefix), so Excel does not
# coerce the assignment into a Double and strip formatting.

$ws.Range('D2').Value = '27.718.18'
$ws.Range('D3').Value = '1.849.28'
$ws.Range('E3').Value = '  -0.90%  '
$ws.Range('D4').Value = "'" + '1.014'
$ws.Range('E4').Value = '  -2.55%  '
$ws.Range('D5').Value = "'" + '319.69'
$ws.Range('E5').Value = '  -1.58%  '
$ws.Range('D6').Value = "'" + '1.011'
$ws.Range('E6').Value = '  -2.47%  '
$ws.Range('D7').Value = "'" + '0.4318'
$ws.Range('E7').Value = '  -2.44%  '
$ws.Range('D8').Value = "'" + '0.3744'
$ws.Range('E8').Value = '  -1.57%  '
$ws.Range('D9').Value = "'" + '0.07355'
$ws.Range('E9').Value = '  -1.60%  '
$ws.Range('D10').Value = "'" + '0.8804'
$ws.Range('E10').Value = '  -0.64%  '
$ws.Range('E11').Value = '  -0.56%  '
$ws.Range('D12').Value = '1.851.38'
$ws.Range('E12').Value = '  -1.18%  '
$ws.Range('D13').Value = "'" + '6.733'
$ws.Range('E13').Value = '  -0.51%  '
$ws.Range('D14').Value = "'" + '5.456'
$ws.Range('E14').Value = '  -1.98%  '
$ws.Range('D15').Value = "'" + '0.07120'
$ws.Range('E15').Value = '  -1.77%  '
$ws.Range('D16').Value = "'" + '88.00'
$ws.Range('E16').Value = '  +4.88%  '
$ws.Range('D17').Value = "'" + '1.015'
$ws.Range('E17').Value = '  -2.62%  '
$ws.Range('D18').Value = "'" + '0.000008995'
$ws.Range('E18').Value = '  -1.82%  '
$ws.Range('D19').Value = "'" + '1.011'
$ws.Range('E19').Value = '  -2.46%  '
$ws.Range('D21').Value = '27.720.89'
$ws.Range('E21').Value = '  -0.27%  '
$ws.Range('E22').Value = '  -1.39%  '
$ws.Range('E23').Value = '  -1.94%  '
$ws.Range('D24').Value = '2.072.55'
$ws.Range('E24').Value = '  -1.50%  '
$ws.Range('D25').Value = "'" + '2.019'
$ws.Range('E25').Value = '  +1.03%  '
$ws.Range('D26').Value = "'" + '155.68'
$ws.Range('E26').Value = '  -2.05%  '
$ws.Range('D27').Value = "'" + '18.62'
$ws.Range('E27').Value = '  -1.41%  '
$ws.Range('D28').Value = "'" + '2.137'
$ws.Range('E28').Value = '  +7.48%  '
$ws.Range('D29').Value = "'" + '5.393'
$ws.Range('E29').Value = '  +1.07%  '
$ws.Range('D30').Value = "'" + '120.57'
$ws.Range('E30').Value = '  +2.24%  '
$ws.Range('D31').Value = "'" + '0.08939'
$ws.Range('E31').Value = '  -1.47%  '
$ws.Range('D32').Value = "'" + '1.229'
$ws.Range('E32').Value = '  +1.20%  '
$ws.Range('D33').Value = "'" + '0.7788'
$ws.Range('E33').Value = '  -0.02%  '
$ws.Range('D34').Value = "'" + '4.566'
$ws.Range('E34').Value = '  -0.17%  '
$ws.Range('D35').Value = "'" + '2.919'
$ws.Range('E35').Value = '  -6.07%  '
$ws.Range('D36').Value = "'" + '1.012'
$ws.Range('E36').Value = '  -2.48%  '
$ws.Range('E37').Value = '  -1.46%  '
$ws.Range('D38').Value = "'" + '0.05335'
$ws.Range('E38').Value = '  -0.34%  '
$ws.Range('E39').Value = '  -1.10%  '
$ws.Range('D40').Value = "'" + '7.252'
$ws.Range('E40').Value = '  +5.14%  '
$ws.Range('D41').Value = "'" + '2.877'
$ws.Range('E41').Value = '  +0.23%  '
$ws.Range('D42').Value = "'" + '0.5160'
$ws.Range('E42').Value = '  -0.89%  '
$ws.Range('D43').Value = "'" + '0.1680'
$ws.Range('E43').Value = '  -0.99%  '
$ws.Range('D44').Value = "'" + '8.915'
$ws.Range('E44').Value = '  +2.61%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = "'" + '10.71'
$ws.Range('E45').Value = '  +0.19%  '
$ws.Range('B46').Value = 'Quant'
$ws.Range('C46').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D46').Value = "'" + '109.53'
$ws.Range('E46').Value = '  -0.37%  '
$ws.Range('D47').Value = "'" + '0.4733'
$ws.Range('E47').Value = '  +0.37%  '
$ws.Range('D48').Value = "'" + '0.06513'
$ws.Range('E48').Value = '  +0.51%  '
$ws.Range('D49').Value = "'" + '1.700'
$ws.Range('E49').Value = '  -1.44%  '
$ws.Range('E50').Value = '  -2.69%  '
$ws.Range('D51').Value = "'" + '1.880'
$ws.Range('E51').Value = '  -2.09%  '
